# "added ifoCAST full series evaluation"
# The naive-QoQ error series got re-evaluated against the full ifoCAST
# history: every data row (row 2..24, columns B..K) now holds a freshly
# recomputed set of errors, and because the matched series gained one
# more leading point, each row's populated values shift left by one
# column relative to before -- so the previous right-most value per row
# is dropped (row 24 ends up with no numeric value left at all).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One CSV line per data row (rows 2-24), values for columns B..K in order.
# A trailing empty field means that cell must be cleared (no longer used).
$dataText = @"
9.881040825240538,-8.249324585355328,-0.0836191634519744,2.07573904723772,-1.358376186195071,-1.416680482701012,0.6261829423615652,0.04427631043098978,0.1402879346210716,-0.1199664748049847
-9.164875883305404,-0.6390580781019446,1.756819337039955,-1.567578570662666,-1.569897914509819,0.5007280727555402,-0.06729608136246609,0.03564018068714139,-0.2211577533692581,0.6161026952046367
-7.804287960688893,-4.017835460506984,-6.209636531594684,-5.320425743141963,-2.555357096058843,-2.587000852527612,-2.071837930259698,-2.012863495286662,-0.9342258126975603,-0.8491485423956758
-2.33765431393776,1.20640823091402,-3.07236942691687,1.426003919130331,-0.6088189506176367,0.3401683158106229,-0.4228531167856981,0.7083467336122384,0.4341740973116463,0.1828529054995925
-2.520798295803975,-1.209498772153247,0.6650949949593515,-0.05858795499538469,0.09291507488738329,-0.1622689761646768,0.6645657133825503,0.5658615947808475,0.2128828273432847,0.3606851936599405
-0.801266318531477,0.9167180388457457,-0.14049400881852,0.1680071242686689,-0.056361815650953,0.7229823595412057,0.6323002786109391,0.2875647859541212,0.4309984855609192,0.5879651212282146
0.8112188269969305,-0.2665992927838778,0.3362551929019006,0.004302276952929285,0.7341151989158385,0.6951630299378696,0.3449163137478644,0.4755036625652983,0.6390389922546331,0.1406021036577823
-0.8710399126316165,0.2720518155472072,0.2958495850101461,0.6964083228299049,0.6865790098644085,0.4286573913452327,0.5103308074933433,0.6656033753069301,0.1861009444043812,0.4585159253925983
0.1844109867082959,0.2434303829804775,0.7550547198149107,0.6915081552015099,0.4175641082963374,0.5230975067184298,0.6743329343542892,0.188602888553233,0.4647918628642795,0.2901314376394514
0.2658583626134861,0.7576582593773593,0.6694404695474946,0.4115821279919825,0.5177542565360187,0.663964400983923,0.1801255664296092,0.4570195545574687,0.2815074641064828,0.589805676679554
0.7790828308345794,0.7896183183467069,0.3223877835781994,0.4880272526874408,0.6725905210138331,0.154185393004174,0.4342409915292412,0.267889428957865,0.5712848794943781,-0.04166000770222517
0.7472976860263175,0.2931833097224661,0.4850066528825197,0.6537292374023527,0.1338542619389697,0.4197108360949934,0.251248383465597,0.553701034065628,-0.05813389123226798,0.6310880987550094
0.6333810503403763,0.5603777146664568,0.4626359510893536,0.1585481745257029,0.4290902396512893,0.196437584854206,0.5348069431769582,-0.07299789204731372,0.6023763881256388,0.3304428394308724
1.011072586656707,0.509746653301712,-0.08173505626510533,0.4650353207171521,0.1943868274146131,0.458238703740928,-0.1003972864225971,0.5760873379279133,0.2870455487175398,
0.8220722402252505,0.051553683470419,0.2859568989009605,0.2225810661851803,0.4942640149333215,-0.1379388566268107,0.5705449988405521,0.2912208776562884,,
0.2871441745782602,0.3706356397752701,0.07698069955146319,0.5042177577925642,-0.1225736869272658,0.5316966638831291,0.2753750686291025,,,
0.6816598262566529,0.1940694739626584,0.3417536163429973,-0.0936147364620642,0.5674395363380327,0.2440474222454754,,,,
0.4408946513667728,0.3594094838808868,-0.1882369755730587,0.5996569245865127,0.2564355480731927,,,,,
0.5990858432970987,-0.1032127321038452,0.4831723462284986,0.2715408197250452,,,,,,
0.0616473449302421,0.4967096184764148,0.2085679007350822,,,,,,,
0.7505586603418228,0.3078859509171186,,,,,,,,
0.3519456421565676,,,,,,,,,
,,,,,,,,,
"@

$rows = $dataText -split "`n" | Where-Object { $_ -ne "" }
$startRow = 2
$startCol = 2  # Column B

for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowIndex = $startRow + $i
    $fields = $rows[$i].Split(",")
    for ($j = 0; $j -lt $fields.Count; $j++) {
        $colIndex = $startCol + $j
        $cellText = $fields[$j]
        if ($cellText -eq "") {
            # Column no longer populated for this row -- remove the cell.
            $ws.Cells.Item($rowIndex, $colIndex).ClearContents()
        } else {
            $ws.Cells.Item($rowIndex, $colIndex).Value = [double]$cellText
        }
    }
}

Write-Host "Done updating ifoCAST matched-error rows"
